$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 282.7879796666667
$ws.Range("H2").Value = 848.363939
$ws.Range("I2").Value = 0.9674521741401267
$ws.Range("J2").Value = 0.9674521741401266
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.454947
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 84219.43979510624
$ws.Range("R2").Value = 757974.9581559561
$ws.Range("S2").Value = 0.8564037967674971
$ws.Range("T2").Value = 0.8564037967674971

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 282.7879796666667
$ws.Range("H3").Value = 848.363939
$ws.Range("I3").Value = 0.9674521741401267
$ws.Range("J3").Value = 0.9674521741401266
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 6883.156798414332
$ws.Range("R3").Value = 61948.41118572899
$ws.Range("S3").Value = 0.06999288561226659
$ws.Range("T3").Value = 0.06999288561226659

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 282.7879796666667
$ws.Range("H4").Value = 848.363939
$ws.Range("I4").Value = 0.9674521741401267
$ws.Range("J4").Value = 0.9674521741401266
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 4037.430158088242
$ws.Range("R4").Value = 36336.87142279417
$ws.Range("S4").Value = 0.04105549176036294
$ws.Range("T4").Value = 0.04105549176036294

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.714696666666668
$ws.Range("H5").Value = 23.14409
$ws.Range("I5").Value = 0.02639291836872237
$ws.Range("J5").Value = 0.02639291836872237
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.454947
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 2297.577967145914
$ws.Range("R5").Value = 20678.20170431323
$ws.Range("S5").Value = 0.02336342415979172
$ws.Range("T5").Value = 0.02336342415979172

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.714696666666668
$ws.Range("H6").Value = 23.14409
$ws.Range("I6").Value = 0.02639291836872237
$ws.Range("J6").Value = 0.02639291836872237
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 187.7783732938856
$ws.Range("R6").Value = 1690.00535964497
$ws.Range("S6").Value = 0.001909465465822921
$ws.Range("T6").Value = 0.001909465465822921

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.714696666666668
$ws.Range("H7").Value = 23.14409
$ws.Range("I7").Value = 0.02639291836872237
$ws.Range("J7").Value = 0.02639291836872237
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 110.14452954902
$ws.Range("R7").Value = 991.3007659411801
$ws.Range("S7").Value = 0.001120028743107736
$ws.Range("T7").Value = 0.001120028743107736

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.799090333333333
$ws.Range("H8").Value = 5.397271
$ws.Range("I8").Value = 0.006154907491150983
$ws.Range("J8").Value = 0.006154907491150983
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.454947
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 535.8020528055151
$ws.Range("R8").Value = 4822.218475249637
$ws.Range("S8").Value = 0.00544842038197843
$ws.Range("T8").Value = 0.005448420381978431

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.799090333333333
$ws.Range("H9").Value = 5.397271
$ws.Range("I9").Value = 0.006154907491150983
$ws.Range("J9").Value = 0.006154907491150983
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 43.79047820010478
$ws.Range("R9").Value = 394.114303800943
$ws.Range("S9").Value = 0.0004452930568532848
$ws.Range("T9").Value = 0.0004452930568532849

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.799090333333333
$ws.Range("H10").Value = 5.397271
$ws.Range("I10").Value = 0.006154907491150983
$ws.Range("J10").Value = 0.006154907491150983
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 25.686033676138
$ws.Range("R10").Value = 231.174303085242
$ws.Range("S10").Value = 0.0002611940523192674
$ws.Range("T10").Value = 0.0002611940523192675
